# Updates cryptos list values (price/volume columns) per the source diff.
# Column D ("Price") values are plain text in the workbook (e.g. "29.308.76" is
# not a valid number), so we force Text format before writing and then restore
# the default "Normal" style so no stray formatting is left on the cell - this
# matches the original inlineStr/shared-string text cells exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.311.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.02%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.873.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.02%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7092"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.48%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.07%  "

# Row 7
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07876"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.18%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3113"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.00%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.18"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.31%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08399"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.30%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.868.40"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.80%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.236"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.16%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7183"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.81%  "

# Row 15
$ws.Range("E15").Value = "  -0.15%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008389"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.42%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.136"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.06%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.318.33"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.00%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.61%  "

# Row 20
$ws.Range("E20").Value = "  +0.04%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.122.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.11%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.05%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.754"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.29%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.05%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1594"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.11%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.23%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.045"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.34%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.08%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.506"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.15%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.411"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.02%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.340"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.19%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.223"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.73%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05354"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.88%  "

# Row 34
$ws.Range("E34").Value = "  +0.85%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.175"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.15%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7470"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.92%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.682"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.04%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01878"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.96%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.245.94"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.02%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.734"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.84%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.503"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.31%  "

# Row 42
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8940"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.56%  "

# Row 43
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "110.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.35%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.92%  "

# Row 45
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000131"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +11.11%  "

# Row 46
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9998"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.00%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.017.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.62%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.799"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.35%  "

# Row 49
$ws.Range("E49").Value = "  -0.04%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.443"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.66%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4348"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.14%  "
